# Applies the cryptos-list refresh described by the commit diff.
# Each D-column numeric-looking value is entered with a leading
# apostrophe (quote-prefix) so Excel stores it as text, matching
# the workbook's existing convention of keeping these as strings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '43.628.26'
$ws.Cells.Item(2, 5).Value = '  +1.13%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.376.32'
$ws.Cells.Item(3, 5).Value = '  +3.28%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.06%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''310.42'
$ws.Cells.Item(5, 5).Value = '  -0.01%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''105.30'
$ws.Cells.Item(6, 5).Value = '  +4.01%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '''0.513'
$ws.Cells.Item(7, 5).Value = '  -4.57%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.01%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  -0.75%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '''36.22'
$ws.Cells.Item(10, 5).Value = '  +0.43%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '''53.48'
$ws.Cells.Item(11, 5).Value = '  +2.61%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  -0.68%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '''0.112'
$ws.Cells.Item(13, 5).Value = '  -0.69%  '

# Row 14
$ws.Cells.Item(14, 5).Value = '  -1.72%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '2.746.19'
$ws.Cells.Item(15, 5).Value = '  +3.59%  '

# Row 16
$ws.Cells.Item(16, 5).Value = '  +4.22%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '2.382.22'
$ws.Cells.Item(17, 5).Value = '  +3.63%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '''0.815'
$ws.Cells.Item(18, 5).Value = '  +0.56%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '43.565.84'
$ws.Cells.Item(19, 5).Value = '  +1.24%  '

# Row 20
$ws.Cells.Item(20, 2).Value = 'Uniswap'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(20, 4).Value = '''6.32'
$ws.Cells.Item(20, 5).Value = '  +3.99%  '

# Row 21
$ws.Cells.Item(21, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(21, 4).Value = '''11.96'
$ws.Cells.Item(21, 5).Value = '  -4.67%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '0.0₃0920'
$ws.Cells.Item(22, 5).Value = '  -0.57%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '''68.44'
$ws.Cells.Item(23, 5).Value = '  -0.04%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '''241.61'
$ws.Cells.Item(24, 5).Value = '  +0.36%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +2.55%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '''2.62'
$ws.Cells.Item(26, 5).Value = '  -0.09%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  -0.06%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  +4.76%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  -2.47%  '

# Row 30
$ws.Cells.Item(30, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(30, 4).Value = '''36.85'
$ws.Cells.Item(30, 5).Value = '  -4.12%  '

# Row 31
$ws.Cells.Item(31, 2).Value = 'Toncoin'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(31, 4).Value = '''2.20'
$ws.Cells.Item(31, 5).Value = '  -5.08%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '''9.57'
$ws.Cells.Item(32, 5).Value = '  -0.72%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '''161.43'
$ws.Cells.Item(33, 5).Value = '  -3.62%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '''5.28'
$ws.Cells.Item(34, 5).Value = '  -0.59%  '

# Row 35
$ws.Cells.Item(35, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(35, 4).Value = '''1.00'
$ws.Cells.Item(35, 5).Value = '  +0.01%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'Celestia'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Cells.Item(36, 4).Value = '''18.36'
$ws.Cells.Item(36, 5).Value = '  +3.25%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  +6.04%  '

# Row 38
$ws.Cells.Item(38, 2).Value = 'RenderToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(38, 4).Value = '''4.71'
$ws.Cells.Item(38, 5).Value = '  +11.49%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(39, 4).Value = '''3.10'
$ws.Cells.Item(39, 5).Value = '  -0.93%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  +0.37%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '''1.94'
$ws.Cells.Item(41, 5).Value = '  +6.26%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  -1.29%  '

# Row 43
$ws.Cells.Item(43, 5).Value = '  -1.69%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '''2.65'
$ws.Cells.Item(44, 5).Value = '  +15.64%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '2.033.69'
$ws.Cells.Item(45, 5).Value = '  +3.03%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '''19.76'
$ws.Cells.Item(46, 5).Value = '  +3.92%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '''0.0290'
$ws.Cells.Item(47, 5).Value = '  +0.42%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  +4.07%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '''10.56'
$ws.Cells.Item(49, 5).Value = '  +7.30%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  +4.46%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  +0.56%  '
